# edit.ps1 -- Applies double line-spacing to the five Discussion-section
# paragraphs, and splits/merges the runs around the rendered-page-break
# marker so <w:lastRenderedPageBreak/> moves from before
# "alignment is still relevant" to before "published papers by Zhang...",
# matching the target diff.
#
# NOTE: paragraphs are located by scanning $d.Paragraphs for a matching
# text prefix (rather than $d.Content.Find, whose resulting Range can
# bleed paragraph-level formatting into the following paragraph when
# later used to set ParagraphFormat/InsertXML on this runtime).

$d = $word.ActiveDocument

function Get-ParagraphByPrefix($doc, $prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.StartsWith($prefix)) {
            return $para
        }
    }
    throw "Could not find paragraph starting with: $prefix"
}

# --- Paragraph 1: "My preliminary tree showed ..." -----------------------
# Double-space the paragraph AND split the run that currently reads
# "...the results matched the ones in the published papers by Zhang..."
# into two runs, moving <w:lastRenderedPageBreak/> onto the second piece.
$p1 = Get-ParagraphByPrefix $d "My preliminary tree showed that SARS-CoV-2"
$p1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="73561082" w14:textId="77777777" w:rsidR="005E1B8E" w:rsidRDefault="00167651"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:t>My preliminary tree showed that SARS-CoV-2 and Pangolin-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CoV</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="005E1B8E"><w:t xml:space="preserve">are closely related to each other. However, the preliminary tree only consisted of the reference genomes. When more taxa were introduced based on specific strains of coronaviruses, the results matched the ones in the </w:t></w:r><w:r w:rsidR="005E1B8E"><w:lastRenderedPageBreak/><w:t>published papers by Zhang et al., and Cohen., 2020. SARS-CoV-2 was found more closely related to Bat-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="005E1B8E"><w:t>CoV</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="005E1B8E"><w:t xml:space="preserve"> (RaTG13) than to Pangolin-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="005E1B8E"><w:t>CoVs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="005E1B8E"><w:t>.</w:t></w:r></w:p>'
$p1.Range.InsertXML($p1Xml)

# --- Paragraph 2: "The Bat-CoV_RaTG13 genome was collected ..." ----------
# Double-space the paragraph AND merge the two runs that used to be split
# by <w:lastRenderedPageBreak/> ("... the " + "alignment is still relevant")
# back into a single run, dropping the now-obsolete page-break marker.
$p2 = Get-ParagraphByPrefix $d "The Bat-CoV_RaTG13 genome was collected"
$p2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="31C8ED16" w14:textId="77777777" w:rsidR="006B142C" w:rsidRDefault="005E1B8E"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:t>The Bat-CoV_RaTG13 genome was collected in 2013 (as mentioned on the NCBI database). Although this strain is relatively older to the SARS-CoV-2 strain that originated in 2019, the alignment is still relevant</w:t></w:r><w:r w:rsidR="006B142C"><w:t xml:space="preserve"> as the RaTG13 genome persisted in the animal pool for over 6 years. Also, the RaTG13 had been around since 2003 and recollected in samples in 2013, indicating that the virus had been survived for at least 16 years. And the more time it survives, the more time it </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="006B142C"><w:t>has to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="006B142C"><w:t xml:space="preserve"> evolve and mutate to infect other hosts. This is very applicable to our knowledge of zoonotic viruses as they evolve and attack human hosts.</w:t></w:r></w:p>'
$p2.Range.InsertXML($p2Xml)

# --- Paragraphs 3-5: just apply double line spacing -----------------------
$p3 = Get-ParagraphByPrefix $d "The phylogenies show a close relationship with RaTG13"
$p3.Range.ParagraphFormat.LineSpacingRule = 2

$p4 = Get-ParagraphByPrefix $d "My phylogenies also showed a relationship between Camel-"
$p4.Range.ParagraphFormat.LineSpacingRule = 2

$p5 = Get-ParagraphByPrefix $d "Overall, my phylogenetic analysis established successful relationships"
$p5.Range.ParagraphFormat.LineSpacingRule = 2

Write-Host "Done."
